# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp text in A1 ---
$ws.Range("A1").Value2 = "Datos actualizados a 22 de Mayo de 2020 a las 11:35"

# --- Row 4 (Estados Unidos) updated totals ---
$ws.Range("B4").Value2 = 1621333
$ws.Range("C4").Value2 = 431
$ws.Range("E4").Value2 = 1142726
$ws.Range("G4").Value2 = 9
$ws.Range("H4").Value2 = 96363

# --- Rows 42-43: Austria overtakes Japon ---
# Row 42 becomes Austria (new figures)
$ws.Range("A42").Value2 = "Austria"
$ws.Range("B42").Value2 = 16436
$ws.Range("C42").Value2 = 32
$ws.Range("D42").Value2 = 15005
$ws.Range("E42").Value2 = 796
$ws.Range("G42").Value2 = 2
$ws.Range("H42").Value2 = 635

# Row 43 becomes Japon (former Austria row's old figures)
$ws.Range("A43").Value2 = "Japon"
$ws.Range("B43").Value2 = 16424
$ws.Range("D43").Value2 = 12672
$ws.Range("E43").Value2 = 2975
$ws.Range("H43").Value2 = 777

# --- Rows 62-65: Oman overtakes Moldavia, Finlandia, Ghana ---
# Row 62 becomes Oman (new figures)
$ws.Range("A62").Value2 = "Oman"
$ws.Range("B62").Value2 = 6794
$ws.Range("C62").Value2 = 424
$ws.Range("D62").Value2 = 1821
$ws.Range("E62").Value2 = 4941
$ws.Range("G62").Value2 = 1
$ws.Range("H62").Value2 = 32

# Row 63 becomes Moldavia (previous row 62 figures)
$ws.Range("A63").Value2 = "Moldavia"
$ws.Range("B63").Value2 = 6704
$ws.Range("D63").Value2 = 2953
$ws.Range("E63").Value2 = 3518
$ws.Range("H63").Value2 = 233

# Row 64 becomes Finlandia (previous row 63 figures)
$ws.Range("A64").Value2 = "Finlandia"
$ws.Range("B64").Value2 = 6493
$ws.Range("C64").Value2 = 0
$ws.Range("D64").Value2 = 4800
$ws.Range("E64").Value2 = 1387
$ws.Range("H64").Value2 = 306

# Row 65 becomes Ghana (previous row 64 figures)
$ws.Range("A65").Value2 = "Ghana"
$ws.Range("B65").Value2 = 6486
$ws.Range("C65").Value2 = 217
$ws.Range("D65").Value2 = 1951
$ws.Range("E65").Value2 = 4504
$ws.Range("G65").Value2 = 0
$ws.Range("H65").Value2 = 31

# --- Row 104 updated totals ---
$ws.Range("B104").Value2 = 1066
$ws.Range("C104").Value2 = 2
$ws.Range("E104").Value2 = 33
